# Update "想去人数" (interest count) values that changed between the two
# publishing runs of the generated workbook (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 40
$ws1.Range("F5").Value = 4945
$ws1.Range("F7").Value = 77
$ws1.Range("F8").Value = 275
$ws1.Range("F9").Value = 39

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 40
$ws4.Range("F9").Value = 4945
$ws4.Range("F11").Value = 77
$ws4.Range("F13").Value = 275
$ws4.Range("F14").Value = 39
